# Applies the "Added email and pdf report" commit:
#  - Inserts two new worksheets, "SendEmail" and "ReportPDF", right after "Login".
#  - Populates both sheets with the SMTP / PDF-report automation settings.
#  - Adds mailto: hyperlinks on the e-mail address cells (Hyperlink style).
#  - Updates the Login sheet's selection / active-sheet bookkeeping.

$wb = $excel.ActiveWorkbook

$loginSheet = $wb.Worksheets.Item("Login")

# --- Create the two new sheets in the right tab order --------------------
$sendEmail = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loginSheet)
$sendEmail.Name = "SendEmail"

$reportPdf = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sendEmail)
$reportPdf.Name = "ReportPDF"

# --- SendEmail sheet -------------------------------------------------------
# Header row (written left-to-right except N1, which was authored before M1 --
# preserved here so shared-string indices line up with the source workbook).
$sendEmail.Range("A1").Value = "Subject"
$sendEmail.Range("B1").Value = "Body"
$sendEmail.Range("C1").Value = "To"
$sendEmail.Range("D1").Value = "From"
$sendEmail.Range("E1").Value = "ServerHostName"
$sendEmail.Range("F1").Value = "ServerPort"
$sendEmail.Range("G1").Value = "UseSSL"
$sendEmail.Range("H1").Value = "Username"
$sendEmail.Range("I1").Value = "Password"
$sendEmail.Range("J1").Value = "SendEmailOnFailure"
$sendEmail.Range("K1").Value = "SendEmailOnSuccess"
$sendEmail.Range("L1").Value = "SendZippedReportOnComplete"
$sendEmail.Range("N1").Value = "PDFReportCustomStyleSheet"

# Data row
$sendEmail.Range("A2").Value = "Build AA Automation 2019"
$sendEmail.Range("B2").WrapText = $true
$sendEmail.Range("C2").Value = "mmargasagayam@abacusnext.com"
$sendEmail.Range("D2").Value = "amicustestmk1@gmail.com"
$sendEmail.Range("E2").Value = "smtp.gmail.com"
$sendEmail.Range("F2").Value = 587
$sendEmail.Range("G2").Value = $true
$sendEmail.Range("H2").Value = "amicustestmk1@gmail.com"
$sendEmail.Range("I2").Value = "0nXTeam123$$"
$sendEmail.Range("J2").Value = $false
$sendEmail.Range("K2").Value = $false
$sendEmail.Range("L2").Value = $false
$sendEmail.Range("M2").Value = $true

# M1 header is written last, matching the source file's shared-string order.
$sendEmail.Range("M1").Value = "SendPDFReportOnComplete"

# Hyperlinks (Excel auto-applies the built-in "Hyperlink" cell style).
$sendEmail.Hyperlinks.Add($sendEmail.Range("C2"), "mailto:mmargasagayam@abacusnext.com")
$sendEmail.Hyperlinks.Add($sendEmail.Range("D2"), "mailto:amicustestmk1@gmail.com")
$sendEmail.Hyperlinks.Add($sendEmail.Range("H2"), "mailto:amicustestmk1@gmail.com")

# Column widths (best effort - approximates the autofit widths from the source).
$sendEmail.Columns.Item(1).ColumnWidth = 24.43
$sendEmail.Columns.Item(2).ColumnWidth = 69.86
$sendEmail.Columns.Item(3).ColumnWidth = 32.57
$sendEmail.Columns.Item(4).ColumnWidth = 25.86
$sendEmail.Columns.Item(5).ColumnWidth = 16.14
$sendEmail.Columns.Item(6).ColumnWidth = 10.43
$sendEmail.Columns.Item(7).ColumnWidth = 7.14
$sendEmail.Columns.Item(8).ColumnWidth = 25.86
$sendEmail.Columns.Item(9).ColumnWidth = 14.29
$sendEmail.Columns.Item(10).ColumnWidth = 19.29
$sendEmail.Columns.Item(11).ColumnWidth = 19.86
$sendEmail.Columns.Item(12).ColumnWidth = 29.43
$sendEmail.Columns.Item(13).ColumnWidth = 26.71
$sendEmail.Columns.Item(14).ColumnWidth = 27.29

$sendEmail.Range("B2").Select()

# --- ReportPDF sheet --------------------------------------------------------
$reportPdf.Range("A1").Value = "PdfName"
$reportPdf.Range("B1").Value = "PdfDirectoryPath"
$reportPdf.Range("C1").Value = "Xml"
$reportPdf.Range("D1").Value = "Details"
$reportPdf.Range("E1").Value = "DeleteRanorexReport"

$reportPdf.Columns.Item(1).ColumnWidth = 9.29
$reportPdf.Columns.Item(2).ColumnWidth = 16.29
$reportPdf.Columns.Item(3).ColumnWidth = 4.43
$reportPdf.Columns.Item(4).ColumnWidth = 7.14
$reportPdf.Columns.Item(5).ColumnWidth = 20.57

$reportPdf.Range("E1").Select()

# --- Login sheet: selection moved, no longer the active tab ---------------
$loginSheet.Range("H5").Select()

# --- Make SendEmail the active tab (matches activeTab=1 / tabSelected) ----
$sendEmail.Range("B2").Select()
$sendEmail.Activate()
